# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Cebolla" (Feria Lagunitas de Puerto Montt)
# right before the current row 257, pushing the existing data (old rows
# 257-285) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 257-258; everything from old row 257 onward
# shifts down to 259 onward (formatting, incl. the date number format on
# column D, is inherited from the row above).
$ws.Rows("257:258").Insert()

# New row 257: Morada(o) / 1a (guarda)
$ws.Cells.Item(257, 1).Value2 = 4
$ws.Cells.Item(257, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(257, 3).Value2 = "Los Lagos"
$ws.Cells.Item(257, 4).Value2 = 44449
$ws.Cells.Item(257, 5).Value2 = 10
$ws.Cells.Item(257, 6).Value2 = 100112004
$ws.Cells.Item(257, 7).Value2 = "Cebolla"
$ws.Cells.Item(257, 8).Value2 = "Morada(o)"
$ws.Cells.Item(257, 9).Value2 = "1a (guarda)"
$ws.Cells.Item(257, 10).Value2 = 150
$ws.Cells.Item(257, 11).Value2 = 14000
$ws.Cells.Item(257, 12).Value2 = 14000
$ws.Cells.Item(257, 13).Value2 = 14000
$ws.Cells.Item(257, 14).Value2 = "`$/malla 18 kilos"
$ws.Cells.Item(257, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(257, 16).Value2 = 778
$ws.Cells.Item(257, 17).Value2 = 18
$ws.Cells.Item(257, 18).Value2 = "Hortaliza"

# New row 258: Sin especificar / 1a (guarda)
$ws.Cells.Item(258, 1).Value2 = 4
$ws.Cells.Item(258, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(258, 3).Value2 = "Los Lagos"
$ws.Cells.Item(258, 4).Value2 = 44449
$ws.Cells.Item(258, 5).Value2 = 10
$ws.Cells.Item(258, 6).Value2 = 100112004
$ws.Cells.Item(258, 7).Value2 = "Cebolla"
$ws.Cells.Item(258, 8).Value2 = "Sin especificar"
$ws.Cells.Item(258, 9).Value2 = "1a (guarda)"
$ws.Cells.Item(258, 10).Value2 = 800
$ws.Cells.Item(258, 11).Value2 = 7000
$ws.Cells.Item(258, 12).Value2 = 7000
$ws.Cells.Item(258, 13).Value2 = 7000
$ws.Cells.Item(258, 14).Value2 = "`$/malla 16 kilos"
$ws.Cells.Item(258, 15).Value2 = "Región de O'Higgins"
$ws.Cells.Item(258, 16).Value2 = 438
$ws.Cells.Item(258, 17).Value2 = 16
$ws.Cells.Item(258, 18).Value2 = "Hortaliza"
